$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.7237827715355806
$ws1.Range("C2").Value = 0.9345454545454546
$ws1.Range("D2").Value = 0.4812734082397004
$ws1.Range("E2").Value = 0.6353522867737948
$ws1.Range("F2").Value = 0.5329738697635836
$ws1.Range("G2").Value = 0.4904220183486239
$ws1.Range("H2").Value = 0.7237827715355806
$ws1.Range("I2").Value = 257
$ws1.Range("J2").Value = 18
$ws1.Range("K2").Value = 516
$ws1.Range("L2").Value = 277

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")
# Row 2 (class 0)
$ws2.Range("B2").Value = 0.6506935687263556
$ws2.Range("C2").Value = 0.9662921348314607
$ws2.Range("D2").Value = 0.7776940467219292
# Row 3 (class 1)
$ws2.Range("B3").Value = 0.9345454545454546
$ws2.Range("C3").Value = 0.4812734082397004
$ws2.Range("D3").Value = 0.6353522867737948
# Row 4 (accuracy)
$ws2.Range("B4").Value = 0.7237827715355806
$ws2.Range("C4").Value = 0.7237827715355806
$ws2.Range("D4").Value = 0.7237827715355806
$ws2.Range("E4").Value = 0.7237827715355806
# Row 5 (macro avg)
$ws2.Range("B5").Value = 0.7926195116359052
$ws2.Range("C5").Value = 0.7237827715355806
$ws2.Range("D5").Value = 0.706523166747862
# Row 6 (weighted avg)
$ws2.Range("B6").Value = 0.7926195116359052
$ws2.Range("C6").Value = 0.7237827715355806
$ws2.Range("D6").Value = 0.7065231667478621

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 516
$ws3.Range("C2").Value = 18
$ws3.Range("B3").Value = 277
$ws3.Range("C3").Value = 257
